# Updates symbol list data per commit "Updated symbol list on Fri Jan 27 08:14:02 UTC 2023 with GitHub Actions"
# Sets each target cell's value as literal text (preserving original text formatting,
# e.g. leading zeros, percent signs, trailing zeros) and restores the default "Normal"
# style so no stray style/number-format indices are introduced.
function Set-CellText($ws, $ref, $text) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws 'D2' '305.78'
Set-CellText $ws 'E2' '-0.31%'
Set-CellText $ws 'G2' '8'
Set-CellText $ws 'D3' '35.66'
Set-CellText $ws 'E3' '-0.72%'
Set-CellText $ws 'G3' '8'
Set-CellText $ws 'D4' '5.040'
Set-CellText $ws 'E4' '-1.15%'
Set-CellText $ws 'G4' '8'
Set-CellText $ws 'D5' '0.07992'
Set-CellText $ws 'E5' '-1.16%'
Set-CellText $ws 'G5' '8'
Set-CellText $ws 'D6' '1.909'
Set-CellText $ws 'E6' '-1.40%'
Set-CellText $ws 'G6' '8'
Set-CellText $ws 'B7' 'KuCoinToken'
Set-CellText $ws 'C7' 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-CellText $ws 'D7' '7.764'
Set-CellText $ws 'E7' '0.24%'
Set-CellText $ws 'G7' '8'
Set-CellText $ws 'B8' 'MXToken'
Set-CellText $ws 'C8' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-CellText $ws 'D8' '0.9210'
Set-CellText $ws 'E8' '-0.73%'
Set-CellText $ws 'G8' '8'
Set-CellText $ws 'B9' 'LiechtensteinCryptoassetsExchange'
Set-CellText $ws 'C9' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-CellText $ws 'D9' '0.1277'
Set-CellText $ws 'E9' '-6.76%'
Set-CellText $ws 'G9' '8'
Set-CellText $ws 'B10' 'WazirX'
Set-CellText $ws 'C10' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-CellText $ws 'D10' '0.1925'
Set-CellText $ws 'E10' '1.43%'
Set-CellText $ws 'G10' '8'
Set-CellText $ws 'B11' 'MandalaExchangeToken'
Set-CellText $ws 'C11' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-CellText $ws 'D11' '0.09122'
Set-CellText $ws 'E11' '-0.89%'
Set-CellText $ws 'G11' '8'
Set-CellText $ws 'B12' 'BitrueCoin'
Set-CellText $ws 'C12' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-CellText $ws 'D12' '0.03444'
Set-CellText $ws 'E12' '0.40%'
Set-CellText $ws 'G12' '8'
Set-CellText $ws 'B13' 'BitMartToken'
Set-CellText $ws 'C13' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-CellText $ws 'D13' '0.09839'
Set-CellText $ws 'E13' '0.00%'
Set-CellText $ws 'G13' '8'
Set-CellText $ws 'B14' 'BitForexToken'
Set-CellText $ws 'C14' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-CellText $ws 'D14' '0.001405'
Set-CellText $ws 'E14' '-2.56%'
Set-CellText $ws 'G14' '8'
Set-CellText $ws 'B15' 'TigerCash'
Set-CellText $ws 'C15' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-CellText $ws 'D15' '0.006159'
Set-CellText $ws 'E15' '6.26%'
Set-CellText $ws 'G15' '8'
Set-CellText $ws 'B16' 'LEO'
Set-CellText $ws 'C16' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-CellText $ws 'D16' '3.719'
Set-CellText $ws 'E16' '2.50%'
Set-CellText $ws 'G16' '8'
Set-CellText $ws 'B17' 'GateToken'
Set-CellText $ws 'C17' 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-CellText $ws 'D17' '4.148'
Set-CellText $ws 'E17' '-1.05%'
Set-CellText $ws 'G17' '8'
Set-CellText $ws 'E18' '12.61%'
Set-CellText $ws 'G18' '8'
Set-CellText $ws 'E19' '-0.28%'
Set-CellText $ws 'G19' '8'
Set-CellText $ws 'D20' '0.1342'
Set-CellText $ws 'E20' '0.62%'
Set-CellText $ws 'G20' '8'
Set-CellText $ws 'D21' '5.162'
Set-CellText $ws 'E21' '5.15%'
Set-CellText $ws 'G21' '8'
Set-CellText $ws 'D22' '0.2595'
Set-CellText $ws 'E22' '6.22%'
Set-CellText $ws 'G22' '8'
Set-CellText $ws 'D23' '0.04436'
Set-CellText $ws 'E23' '-0.12%'
Set-CellText $ws 'G23' '8'
Set-CellText $ws 'D24' '0.001234'
Set-CellText $ws 'E24' '0.93%'
Set-CellText $ws 'G24' '8'
Set-CellText $ws 'D25' '0.004635'
Set-CellText $ws 'E25' '-3.81%'
Set-CellText $ws 'G25' '8'
Set-CellText $ws 'E26' '-4.06%'
Set-CellText $ws 'G26' '8'
Set-CellText $ws 'D27' '0.0004443'
Set-CellText $ws 'E27' '41.71%'
Set-CellText $ws 'G27' '8'
Set-CellText $ws 'G28' '8'
Set-CellText $ws 'G29' '8'
Set-CellText $ws 'G30' '8'
Set-CellText $ws 'G31' '8'
Set-CellText $ws 'G32' '8'
Set-CellText $ws 'G33' '8'
Set-CellText $ws 'G34' '8'
Set-CellText $ws 'G35' '8'
Set-CellText $ws 'G36' '8'
Set-CellText $ws 'G37' '8'
Set-CellText $ws 'G38' '8'
Set-CellText $ws 'D39' '0.01943'
Set-CellText $ws 'G39' '8'
Set-CellText $ws 'D40' '0.05326'
Set-CellText $ws 'E40' '8.30%'
Set-CellText $ws 'G40' '8'
Set-CellText $ws 'D41' '0.007621'
Set-CellText $ws 'E41' '-0.53%'
Set-CellText $ws 'G41' '8'
Set-CellText $ws 'D42' '0.01020'
Set-CellText $ws 'E42' '-1.53%'
Set-CellText $ws 'G42' '8'
Set-CellText $ws 'D43' '0.1356'
Set-CellText $ws 'E43' '-1.61%'
Set-CellText $ws 'G43' '8'
Set-CellText $ws 'E44' '2.15%'
Set-CellText $ws 'G44' '8'
Set-CellText $ws 'D45' '0.009899'
Set-CellText $ws 'E45' '-10.72%'
Set-CellText $ws 'G45' '8'
Set-CellText $ws 'D46' '0.00006134'
Set-CellText $ws 'E46' '-4.45%'
Set-CellText $ws 'G46' '8'
Set-CellText $ws 'E47' '-0.22%'
Set-CellText $ws 'G47' '8'
Set-CellText $ws 'D48' '65.22'
Set-CellText $ws 'E48' '2.60%'
Set-CellText $ws 'G48' '8'
Set-CellText $ws 'D49' '0.001659'
Set-CellText $ws 'E49' '39.03%'
Set-CellText $ws 'G49' '8'
Set-CellText $ws 'D50' '0.00002100'
Set-CellText $ws 'E50' '-0.22%'
Set-CellText $ws 'G50' '8'
Set-CellText $ws 'D51' '0.0002000'
Set-CellText $ws 'E51' '-0.22%'
Set-CellText $ws 'G51' '8'

Write-Host "Applied 144 cell updates"
